# The deck's design theme (ppt/theme/theme1.xml, linked from the
# presentation's one SlideMaster) is re-coloured from the "Integral"
# palette to the stock "Office" palette — i.e. it is swapped for the
# theme that used to live in ppt/theme/theme2.xml (the font scheme and
# format/effect scheme are already identical between the two theme
# parts, so only the twelve colour-scheme slots actually change).
#
# PowerPoint's object model exposes those twelve theme colours through
# ThemeColorScheme (Dark1, Light1, Dark2, Light2, Accent1-6, Hyperlink,
# FollowedHyperlink, in that fixed order), settable via .RGB on each
# item. We drive that through the Slide's ThemeColorScheme, which maps
# straight onto the presentation's theme part.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# PowerPoint's RGB() long is little-endian (0x00BBGGRR), so R + G*256 + B*65536.
$tcs.Item(1).RGB  = 0x000000   # Dark 1   -> 000000
$tcs.Item(2).RGB  = 0xFFFFFF   # Light 1  -> FFFFFF
$tcs.Item(3).RGB  = 0x6A5444   # Dark 2   -> 44546A
$tcs.Item(4).RGB  = 0xE6E6E7   # Light 2  -> E7E6E6
$tcs.Item(5).RGB  = 0xD59B5B   # Accent 1 -> 5B9BD5
$tcs.Item(6).RGB  = 0x317DED   # Accent 2 -> ED7D31
$tcs.Item(7).RGB  = 0xA5A5A5   # Accent 3 -> A5A5A5
$tcs.Item(8).RGB  = 0x00C0FF   # Accent 4 -> FFC000
$tcs.Item(9).RGB  = 0xC47244   # Accent 5 -> 4472C4
$tcs.Item(10).RGB = 0x47AD70   # Accent 6 -> 70AD47
$tcs.Item(11).RGB = 0xC16305   # Hyperlink          -> 0563C1
$tcs.Item(12).RGB = 0x724F95   # Followed Hyperlink -> 954F72
